# Append the 24 Feb 2021 hangboard training entry as a new row (row 38)
# at the bottom of the log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "24 Feb 2021"
$ws.Range("B38").Value = "0,6"
$ws.Range("C38").Value = "-12.5,6"
$ws.Range("D38").Value = "-22.5,6"
$ws.Range("E38").Value = "-42.5,5,8"
$ws.Range("F38").Value = "-30,6"
$ws.Range("G38").Value = "-15,5,9"
$ws.Range("H38").Value = "-40,4,8,5"
$ws.Range("I38").Value = "-22.5,5,8"
$ws.Range("J38").Value = "-30,5,7"

# Move the active selection to the newly added last cell, matching where
# the author's cursor ended up after typing the new row.
$ws.Range("J38").Select()
